# Update the CDA Logical model metadata for ST.r2b
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row after "Contact" (row 10) for the "Jurisdiction" property, pushing
# Description/Purpose/... etc down by one row.
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Bump version + date metadata values.
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"
